$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price cells remain stored as text (matching original inlineStr type)
$textCells = @("D5", "D6", "D7", "D8", "D10", "D16", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D29", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell content updates as described by the diff
$ws.Range("D2").Value = "58.005.07"
$ws.Range("E2").Value = "  +3.62%  "
$ws.Range("D3").Value = "2.458.20"
$ws.Range("E3").Value = "  +3.26%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "159.45"
$ws.Range("E5").Value = "  +7.90%  "
$ws.Range("D6").Value = "494.02"
$ws.Range("E6").Value = "  +3.75%  "
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").Value = "0.611"
$ws.Range("E8").Value = "  +22.95%  "
$ws.Range("D9").Value = "2.486.76"
$ws.Range("E9").Value = "  +4.42%  "
$ws.Range("D10").Value = "6.31"
$ws.Range("E10").Value = "  +14.81%  "
$ws.Range("E11").Value = "  +5.20%  "
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D14").Value = "2.876.11"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("D15").Value = "57.879.42"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").Value = "21.16"
$ws.Range("E16").Value = "  +4.44%  "
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").Value = "2.476.50"
$ws.Range("E18").Value = "  +4.23%  "
$ws.Range("E19").Value = "  +6.72%  "
$ws.Range("D20").Value = "327.21"
$ws.Range("E20").Value = "  +4.84%  "
$ws.Range("D21").Value = "10.17"
$ws.Range("E21").Value = "  +4.74%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.97"
$ws.Range("E23").Value = "  +5.58%  "
$ws.Range("D24").Value = "58.64"
$ws.Range("E24").Value = "  +3.22%  "
$ws.Range("E25").Value = "  +3.49%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.163"
$ws.Range("E26").Value = "  +3.68%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").Value = "2.547.16"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "7.42"
$ws.Range("E29").Value = "  +2.17%  "
$ws.Range("D30").Value = "0.0₃0807"
$ws.Range("E30").Value = "  +4.88%  "
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("E32").Value = "  +5.63%  "
$ws.Range("D33").Value = "150.95"
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("E35").Value = "  +8.84%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "3.83"
$ws.Range("E36").Value = "  +6.71%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.16"
$ws.Range("E37").Value = "  +5.26%  "
$ws.Range("D38").Value = "0.842"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "1.42"
$ws.Range("E39").Value = "  +5.05%  "
$ws.Range("D40").Value = "3.60"
$ws.Range("E40").Value = "  +6.79%  "
$ws.Range("D41").Value = "34.41"
$ws.Range("E41").Value = "  +2.99%  "
$ws.Range("E42").Value = "  +7.18%  "
$ws.Range("D43").Value = "281.39"
$ws.Range("E43").Value = "  +10.84%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "0.609"
$ws.Range("E44").Value = "  +4.28%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "0.992"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").Value = "0.0543"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "4.76"
$ws.Range("E47").Value = "  +5.16%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0232"
$ws.Range("E48").Value = "  +4.59%  "
$ws.Range("D49").Value = "10.25"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("E50").Value = "  +6.94%  "
$ws.Range("D51").Value = "0.689"
$ws.Range("E51").Value = "  +11.17%  "
